# Edit script: update Paysheet schedule rows to reflect new class ordering
# (a new "4/5" date group is inserted, and three more weeks (through 4/29)
# now have their own "ext" closing entries, pushing the sheet from
# A1:I45 to A1:I48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Left block data (columns A:C) for rows 11-48
$leftData = @(
    @("4/3","ext",".5"),
    @("4/5","w67","1"),
    @("4/5","w46","1"),
    @("4/5","ext",".5"),
    @("4/8","w75","1"),
    @("4/8","w45","1"),
    @("4/8","ext",".5"),
    @("4/10","w61","1"),
    @("4/10","w71","1"),
    @("4/10","pe",".5"),
    @("4/10","ext",".5"),
    @("4/12","f70","2"),
    @("4/12","w52","1"),
    @("4/12","ext",".5"),
    @("4/15","w57","1"),
    @("4/15","w48","1"),
    @("4/15","ext",".5"),
    @("4/17","w60","1"),
    @("4/17","pe",".5"),
    @("4/17","w66","1"),
    @("4/17","w49","2"),
    @("4/17","ext",".5"),
    @("4/19","w67","1"),
    @("4/19","w46","1"),
    @("4/19","ext",".5"),
    @("4/22","w75","1"),
    @("4/22","w45","1"),
    @("4/22","ext",".5"),
    @("4/24","w61","1"),
    @("4/24","w71","1"),
    @("4/24","pe",".5"),
    @("4/24","ext",".5"),
    @("4/26","f70","2"),
    @("4/26","w52","1"),
    @("4/26","ext",".5"),
    @("4/29","w57","1"),
    @("4/29","w48","1"),
    @("4/29","ext",".5")
)

# Right block data (columns F:H) for rows 11-48
$rightData = @(
    @("4/5","f70","2"),
    @("4/5","w52","1"),
    @("4/5","ext",".5"),
    @("4/8","w57","1"),
    @("4/8","w48","1"),
    @("4/8","ext",".5"),
    @("4/10","w60","1"),
    @("4/10","pe",".5"),
    @("4/10","w66","1"),
    @("4/10","w49","2"),
    @("4/10","ext",".5"),
    @("4/12","w67","1"),
    @("4/12","w46","1"),
    @("4/12","ext",".5"),
    @("4/15","w75","1"),
    @("4/15","w45","1"),
    @("4/15","ext",".5"),
    @("4/17","w61","1"),
    @("4/17","w71","1"),
    @("4/17","pe",".5"),
    @("4/17","ext",".5"),
    @("4/19","f70","2"),
    @("4/19","w52","1"),
    @("4/19","ext",".5"),
    @("4/22","w57","1"),
    @("4/22","w48","1"),
    @("4/22","ext",".5"),
    @("4/24","w60","1"),
    @("4/24","pe",".5"),
    @("4/24","w66","1"),
    @("4/24","w49","2"),
    @("4/24","ext",".5"),
    @("4/26","w67","1"),
    @("4/26","w46","1"),
    @("4/26","ext",".5"),
    @("4/29","w75","1"),
    @("4/29","w45","1"),
    @("4/29","ext",".5")
)

$startRow = 11
$endRow = $startRow + $leftData.Count - 1

# The "Length" columns (C and H) hold values that look numeric ("1", "2",
# ".5"), but in this workbook they are stored as text, just like every
# other cell. Mark the ranges as Text before writing so Excel does not
# silently convert them to real numbers.
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"
$ws.Range("H$startRow`:H$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $leftData.Count; $i++) {
    $r = $startRow + $i
    $row = $leftData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

for ($i = 0; $i -lt $rightData.Count; $i++) {
    $r = $startRow + $i
    $row = $rightData[$i]
    $ws.Cells.Item($r, 6).Value = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
}
